$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -5686.582309231221
$ws.Range("C2").Value = 13330.70863251111
$ws.Range("D2").Value = -6993.55224781001
$ws.Range("E2").Value = -650.5740754699494
$ws.Range("F2").Value = 60.83066002748302
$ws.Range("G2").Value = 45.00241456370861
$ws.Range("H2").Value = 61.37293966472015
$ws.Range("I2").Value = 44.69642864725432
$ws.Range("J2").Value = 60.92981297429748
$ws.Range("K2").Value = 45.03907235094567
$ws.Range("L2").Value = 53.12368901318746
$ws.Range("O2").Value = 45.95613052498413
$ws.Range("P2").Value = 54.99078892544287
$ws.Range("R2").Value = 5.132400196711578
$ws.Range("S2").Value = -11.41957658593611
$ws.Range("T2").Value = 6.287176389224532
$ws.Range("X2").Value = -53.75641413331196
$ws.Range("Y2").Value = -70.18623935424107
$ws.Range("Z2").Value = -45.53134162955997
$ws.Range("AE2").Value = -5.132400196711578
$ws.Range("AF2").Value = 6.287176389224532
$ws.Range("AG2").Value = 5.132400196711578
$ws.Range("AH2").Value = -11.41957658593611
$ws.Range("AI2").Value = 6.287176389224532
$ws.Range("AJ2").Value = 5.132400196711578
$ws.Range("AK2").Value = -6.287176389224532
$ws.Range("AL2").Value = 8.214912610464552
$ws.Range("AM2").Value = -12.32744886234054
$ws.Range("AN2").Value = -53.75641413331196
$ws.Range("AO2").Value = -70.18623935424107
$ws.Range("AP2").Value = -45.53134162955997
$ws.Range("AQ2").Value = -8.214912610464552
$ws.Range("AR2").Value = 12.32744886234054
$ws.Range("AS2").Value = 60.83066002748302
$ws.Range("AT2").Value = 60.83066002748302
$ws.Range("AU2").Value = 61.37293966472015
$ws.Range("AV2").Value = 61.37293966472015
$ws.Range("AW2").Value = 61.37293966472015
$ws.Range("AX2").Value = 60.92981297429748
$ws.Range("AY2").Value = 60.92981297429748
$ws.Range("AZ2").Value = 45.00241456370861
$ws.Range("BA2").Value = 45.00241456370861
$ws.Range("BB2").Value = 44.69642864725432
$ws.Range("BC2").Value = 44.63993792125859
$ws.Range("BD2").Value = 44.74254362721706
$ws.Range("BE2").Value = 45.03907235094567
$ws.Range("BF2").Value = 45.03907235094567
$ws.Range("BG2").Value = 53.12368901318746
$ws.Range("BJ2").Value = 45.95613052498413
$ws.Range("BK2").Value = 54.99078892544287
